# Loan RBI, Variable Instalments
# Insert a new (blank) column before column M on the "Repayment schedule"
# sheet, shifting the existing "In Advance"/"Late"/heading/heading/
# "Outstanding" columns one to the right, and update the repayment-schedule
# figures to reflect the variable-instalment split (extra principal amount
# moved into the new column).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at M (pushes old M..Q to N..R).
$ws.Range("M1").EntireColumn.Insert()

# Row 3 value updates after the column insert.
# Interest/Fees due (K3) drops from 10045.16 to 45.16 ...
$ws.Range("K3").Value = 45.16
# ... because the 10000 principal portion now shows up separately in the
# newly-shifted "In Advance" column (N3, old M3).
$ws.Range("N3").Value = 10000

# Make "Repayment schedule" the active/selected sheet & cell (was
# "NewLoanInput" before).
$ws.Activate() | Out-Null
$ws.Range("G13").Select() | Out-Null
